$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting so these remain plain text (shared string) values
# instead of being auto-converted to a number / date by Excel.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"

# Update row 2 values: Mobile Number, Transaction Date, Transaction Amount, Milk Quantity (L)
$ws.Range("B2").Value = "9420208901"
$ws.Range("C2").Value = "2025-03-26"
$ws.Range("D2").Value = 780
$ws.Range("E2").Value = 150

# Restore the default "Normal" style on these two cells so they keep the
# same (unstyled) appearance as the rest of the row, now that the text
# value has been safely stored.
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Style = "Normal"

# Delete row 3 entirely (shifts rows up, removing the extra data row)
$ws.Rows.Item(3).Delete()
